# 案件情報.xlsx - append latest scrape run (2026-01-13 01:23:58 JST)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# --- 1. Insert a new row at position 12 (the sheet stays sorted by priority
#        score descending; the new item's score of 75 slots in between the
#        existing 85 and 18 rows), pushing the old rows 12-13 down to 13-14.
$ws.Rows.Item(12).Insert()

# --- 2. Refresh the "取得日時" timestamp on every data row (2-14) to the
#        new scrape time.
$timestamp = "2026-01-13 01:23:58"
for ($r = 2; $r -le 14; $r++) {
    $ws.Range("A$r").Value = $timestamp
}

# --- 3. Row 5 (Shopee API tool) got re-priced and re-scored on this run.
$ws.Range("D5").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("G5").Value = 305

# --- 4. Populate the brand-new row 12 with the newly scraped listing.
$ws.Range("B12").Value = "GoogleCloudを利用したアジャイル開発共通基盤のSREエンジニアの募集"
$ws.Range("C12").Value = "システム開発"
$ws.Range("D12").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E12").Value = "期限情報なし"
$ws.Range("F12").Value = "https://www.lancers.jp/work/detail/5457458"
$ws.Range("G12").Value = 75
$ws.Range("H12").Value = "◆開発"

# --- 5. The row-insert above correctly shifts cell contents/styles down,
#        but this engine does NOT shift the worksheet's Hyperlinks
#        collection along with it, which would leave the F12/F13 link
#        objects stale (still bound to F12/F13 even though that data now
#        lives at F13/F14). Rebuild the whole Hyperlinks collection from
#        scratch, in the correct final row order, to avoid stale/duplicate
#        link relationships.
$ws.Hyperlinks.Delete()

$urls = @{
    2  = "https://www.lancers.jp/work/detail/5423720"
    3  = "https://www.lancers.jp/work/detail/5434128"
    4  = "https://www.lancers.jp/work/detail/5427956"
    5  = "https://www.lancers.jp/work/detail/5469483"
    6  = "https://www.lancers.jp/work/detail/5439158"
    7  = "https://www.lancers.jp/work/detail/5469379"
    8  = "https://www.lancers.jp/work/detail/5469627"
    9  = "https://www.lancers.jp/work/detail/5469430"
    10 = "https://www.lancers.jp/work/detail/5454210"
    11 = "https://www.lancers.jp/work/detail/5469522"
    12 = "https://www.lancers.jp/work/detail/5457458"
    13 = "https://www.lancers.jp/work/detail/5469298"
    14 = "https://www.lancers.jp/work/detail/5469531"
}
for ($r = 2; $r -le 14; $r++) {
    $cell = $ws.Range("F$r")
    $cell.Value = $urls[$r]
    $ws.Hyperlinks.Add($cell, $urls[$r]) | Out-Null
    $cell.Style = "Hyperlink"
}

# --- 6. Column D ("価格") got a bit wider to fit the new, longer price text.
$ws.Columns.Item(4).ColumnWidth = 29.17
